$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New weekly data: insert two rows at the top of the data block (row 62),
# pushing all existing data rows down by two (62-89 -> 64-91).
$ws.Rows("62:63").Insert()

# Row 62 - Lapins / Especial
$ws.Cells.Item(62, 1).Value = 7
$ws.Cells.Item(62, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(62, 3).Value = "Ñuble"
$ws.Cells.Item(62, 4).Value = 44572
$ws.Cells.Item(62, 5).Value = 16
$ws.Cells.Item(62, 6).Value = "Fruta"
$ws.Cells.Item(62, 7).Value = 100103
$ws.Cells.Item(62, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(62, 9).Value = 100103001
$ws.Cells.Item(62, 10).Value = "Cereza"
$ws.Cells.Item(62, 11).Value = "Lapins"
$ws.Cells.Item(62, 12).Value = "Especial"
$ws.Cells.Item(62, 13).Value = 100
$ws.Cells.Item(62, 14).Value = 7000
$ws.Cells.Item(62, 15).Value = 7500
$ws.Cells.Item(62, 16).Value = 7250
$ws.Cells.Item(62, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(62, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(62, 19).Value = 725
$ws.Cells.Item(62, 20).Value = 10

# Row 63 - Lapins / Primera
$ws.Cells.Item(63, 1).Value = 7
$ws.Cells.Item(63, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(63, 3).Value = "Ñuble"
$ws.Cells.Item(63, 4).Value = 44572
$ws.Cells.Item(63, 5).Value = 16
$ws.Cells.Item(63, 6).Value = "Fruta"
$ws.Cells.Item(63, 7).Value = 100103
$ws.Cells.Item(63, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(63, 9).Value = 100103001
$ws.Cells.Item(63, 10).Value = "Cereza"
$ws.Cells.Item(63, 11).Value = "Lapins"
$ws.Cells.Item(63, 12).Value = "Primera"
$ws.Cells.Item(63, 13).Value = 160
$ws.Cells.Item(63, 14).Value = 6000
$ws.Cells.Item(63, 15).Value = 6500
$ws.Cells.Item(63, 16).Value = 6250
$ws.Cells.Item(63, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(63, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(63, 19).Value = 625
$ws.Cells.Item(63, 20).Value = 10
